$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = 3507
$ws.Range("E3").Value  = 606
$ws.Range("E4").Value  = 325
$ws.Range("E5").Value  = 3034
$ws.Range("E6").Value  = 4485
$ws.Range("E7").Value  = 3488
$ws.Range("E8").Value  = 10304
$ws.Range("E9").Value  = 14669
$ws.Range("E10").Value = 2397
$ws.Range("E11").Value = 6528
$ws.Range("E12").Value = 3064
$ws.Range("E13").Value = 8665
$ws.Range("E14").Value = 4663
